$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns, matching the style of the existing header row
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill season record (Wins/Losses/Ties) for every player row
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 29).Value = 59
    $ws.Cells.Item($row, 30).Value = 103
    $ws.Cells.Item($row, 31).Value = 0
}
